$d = $word.ActiveDocument

# 1. Merge "Coraline" + " e o mundo secreto" into a single run "Coraline e o mundo secreto"
#    and drop the spell-check proofErr markers that wrapped "Coraline".
#    Deleting the whole paragraph (including its proofErr markup) and reinserting a clean
#    paragraph avoids leaving stray <w:proofErr/> elements behind.
$coralineRange = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Coraline" -or $t -eq "Coraline e o mundo secreto") {
        $coralineRange = $p.Range
    }
}
$coralineRange.Delete()

$nextPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "O estranho mundo de Jack") {
        $nextPara = $p
    }
}
$nextRange = $nextPara.Range
$nextRange.Collapse(1)
$nextRange.InsertBefore("Coraline e o mundo secreto`r")

# 2. Add a new paragraph "Superman" right after "A fuga das galinhas"
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "A fuga das galinhas") {
        $target = $p
    }
}
$r = $target.Range
$r.InsertParagraphAfter()
$newRange = $d.Range($r.End, $r.End)
$newRange.Text = "Superman"
